# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to stay a TEXT cell
# (the sheet stores Price/Volume columns as text, and some values such as
# "228.25" would otherwise be auto-coerced to a number by the COM Value
# setter). We briefly force a text number-format so the value is stored
# as a string, then clear the format again so no stray style index is
# left behind on the cell (matching the original formatting).
function Set-CellText {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.ClearFormats()
}

# --- Price (D) / Volume(1h) (E) updates ---
Set-CellText "D2" '39.422.88'
Set-CellText "E2" '  +1.60%  '
Set-CellText "D3" '2.159.68'
Set-CellText "E3" '  +3.14%  '
Set-CellText "E4" '  +0.06%  '
Set-CellText "D5" '228.25'
Set-CellText "E5" '  -0.41%  '
Set-CellText "E6" '  +0.95%  '
Set-CellText "D7" '64.33'
Set-CellText "E7" '  +4.94%  '
Set-CellText "E8" '  +0.07%  '
Set-CellText "E9" '  +3.02%  '
Set-CellText "E10" '  +2.02%  '
Set-CellText "E11" '  -0.10%  '
Set-CellText "E12" '  +4.14%  '
Set-CellText "D13" '2.480.89'
Set-CellText "E13" '  +3.20%  '
Set-CellText "D14" '22.32'
Set-CellText "E14" '  +1.14%  '
Set-CellText "E15" '  +0.95%  '
Set-CellText "E16" '  +1.20%  '
Set-CellText "D17" '2.154.64'
Set-CellText "E17" '  +3.46%  '
Set-CellText "D18" '39.350.80'
Set-CellText "E18" '  +1.58%  '
Set-CellText "E19" '  +0.07%  '
Set-CellText "E20" '  +0.60%  '
Set-CellText "D21" '0.0₃0856'
Set-CellText "E21" '  +1.75%  '
Set-CellText "D22" '231.59'
Set-CellText "E22" '  +1.69%  '
Set-CellText "E23" '  +0.11%  '
Set-CellText "E24" '  +5.95%  '
Set-CellText "D25" '2.36'
Set-CellText "E25" '  +0.61%  '
Set-CellText "D26" '172.29'
Set-CellText "E26" '  +0.46%  '
Set-CellText "D27" '9.52'
Set-CellText "E27" '  -0.15%  '
Set-CellText "E28" '  +1.05%  '
Set-CellText "D29" '19.94'
Set-CellText "E29" '  +3.03%  '
Set-CellText "E30" '  -1.33%  '
Set-CellText "D31" '2.67'
Set-CellText "E31" '  +8.59%  '
Set-CellText "E32" '  +1.14%  '
Set-CellText "D33" '4.62'
Set-CellText "E33" '  +2.21%  '
Set-CellText "E34" '  -0.16%  '
Set-CellText "D35" '7.06'
Set-CellText "E35" '  +8.97%  '
Set-CellText "E36" '  +1.15%  '
Set-CellText "E37" '  +0.17%  '
Set-CellText "E38" '  +0.17%  '
Set-CellText "E39" '  +0.12%  '
Set-CellText "E40" '  +1.31%  '
Set-CellText "D41" '103.71'
Set-CellText "E41" '  +2.76%  '
Set-CellText "E42" '  -1.07%  '
Set-CellText "D43" '1.540.26'
Set-CellText "E43" '  +0.29%  '
Set-CellText "E44" '  +4.23%  '
Set-CellText "E49" '  +5.48%  '
Set-CellText "D50" '2.363.91'
Set-CellText "E50" '  +3.19%  '
Set-CellText "D51" '2.96'
Set-CellText "E51" '  -0.30%  '

# --- Rows 45-48 were re-ranked: each coin moved up one slot and FraxShare
#     (previously row 45) wrapped around to row 48. Update Coin name,
#     Link, Price and Volume(1h) for each of these rows. ---
Set-CellText "B45" 'FTXToken'
Set-CellText "C45" 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-CellText "D45" '4.31'
Set-CellText "E45" '  +4.27%  '

Set-CellText "B46" 'HuobiToken'
Set-CellText "C46" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-CellText "D46" '2.82'
Set-CellText "E46" '  +0.67%  '

Set-CellText "B47" 'Cronos'
Set-CellText "C47" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-CellText "D47" '0.0925'
Set-CellText "E47" '  +1.38%  '

Set-CellText "B48" 'FraxShare'
Set-CellText "C48" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-CellText "D48" '7.82'
Set-CellText "E48" '  +1.80%  '
